$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark column H (Absent) as 1 for rows 3 through 18
for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Row 3 also has the Invalid column (G) set to 1
$ws.Cells.Item(3, 7).Value = 1

# Row 12 is special: Total Attendance Count (D) and Real (E) are set to 1,
# and Absent (H) stays at 0 (overriding the general loop above).
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 8).Value = 0
